# Add a sequential numeric "ID" column (A2, A3, ...) for all data rows,
# matching the existing "ID" header already present in A1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}
